# verso-3a.pptx : shrink the "[conteudo]" placeholder text in the
# methodology table from 12pt down to 10pt (the literal brackets were
# already 12pt while the "conteudo" run in between was already 10pt).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the table shape ("Tabela 2") rather than assuming a fixed index.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

$tbl = $tableShape.Table
$cell = $tbl.Cell(2, 1)
$tr = $cell.Shape.TextFrame.TextRange

# Only the runs still at 12pt ("[" and "]") need to drop to 10pt;
# the "conteudo" run in between is left untouched since it is 10pt already.
$runs = $tr.Runs()
for ($i = 1; $i -le $runs.Count; $i++) {
    $run = $runs.Item($i)
    if ($run.Font.Size -eq 12) {
        $run.Font.Size = 10
    }
}
